$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the newly filled-in car-detail columns (E:T) for rows that were
# previously missing them. Keys are row numbers, values are the ordered
# list of values for columns E,F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T.
$rowsData = @{
    127 = @(4,    8,  480, 2, 2, 3,   300, 500, 1100, 0,   120, 478, 195, 120, 275, 6)
    129 = @(5.2, 10,  560, 2, 2, 3.7, 325, 540, 1520, 110, 90,  435, 190, 117, 256, 6)
    130 = @(2,    4,  280, 5, 4, 5.9, 250, 383, 1320, 430, 55,  449, 177, 145, 263, 6)
    131 = @(4,    8,  510, 2, 2, 3.6, 314, 685, 1530, 350, 73,  447, 194, 127, 270, 8)
    132 = @(4.2,  8,  298, 4, 2, 6.2, 250, 411, 1645, 300, 71,  479, 189, 132, 275, 6)
    133 = @(5,    8,  412, 4, 2, 4.5, 250, 529, 1635, 379, 61,  478, 188, 142, 272, 6)
    134 = @(4.2,  8,  650, 2, 2, 3,   370, 850, 1468, 100, 120, 446, 200, 111, 270, 6)
    135 = @(5.5, 12,  442, 4, 2, 5.2, 302, 550, 1690, 316, 110, 473, 192, 130, 260, 6)
    136 = @(2,    4,  230, 4, 2, 5.3, 250, 370, 1335, 305, 55,  418, 183, 135, 251, 6)
    137 = @(5.5,  8,  491, 2, 2, 3.2, 295, 658, 1110, 0,   85,  450, 205, 115, 271, 6)
    138 = @(4.5999999999999996, 8, 232, 6, 4, 7.9, 200, 389, 1905, 583, 72, 539, 199, 144, 291, 4)
    139 = @(2,    4,  150, 5, 5, 8.6999999999999993, 217, 200, 1260, 330, 50, 423, 175, 143, 266, 6)
    140 = @(4.9000000000000004, 8, 400, 2, 2, 4.7, 250, 500, 1660, 203, 73, 440, 183, 132, 251, 6)
    141 = @(4.2,  8,  298, 4, 2, 6.6, 250, 411, 1775, 307, 75,  476, 183, 131, 259, 6)
    142 = @(3,    6,  220, 5, 4, 6.9, 245, 300, 1425, 445, 70,  455, 177, 143, 265, 5)
}

$columns = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($row in $rowsData.Keys) {
    $values = $rowsData[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range($columns[$i] + $row).Value = $values[$i]
    }
}

# Restore the view: move the active selection to J16 instead of F130
# (this also drops the scrolled-down "topLeftCell=A100" view state).
$ws.Range("J16").Select()
